$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header values (B1:E1) ---
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# --- Row 2 (CON) ---
# C2 is deleted entirely (its column "moves" into B2's slot conceptually,
# but in the sheet XML C2 just disappears)
$ws.Range("B2").Value = 136.81760523438663
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 134.85450354383354
$ws.Range("E2").Value = 178.18826857376496

# --- Row 3 (STR) ---
# B3 is deleted entirely
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 176.56955643873493
$ws.Range("D3").Value = 132.33390606648163
$ws.Range("E3").Value = 173.33385639635713

# --- Selection now only covers B1:E3 instead of the whole used range ---
$ws.Range("B1:E3").Select()

Write-Output "Applied passive tweak edits"
